$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (prices, volume %, swapped KuCoinToken/FTXToken rows)
# Values are written as text (leading apostrophe) to match the source sheet,
# which stores these columns as inline strings, then the style is reset to
# "Normal" so Excel's auto-applied @ (Text) number format does not stick.
$updates = @(
    @{ Cell = 'D2'; Value = '310.46' }
    @{ Cell = 'E2'; Value = '1.61%' }
    @{ Cell = 'D3'; Value = '37.28' }
    @{ Cell = 'E3'; Value = '0.55%' }
    @{ Cell = 'D4'; Value = '5.116' }
    @{ Cell = 'E4'; Value = '0.45%' }
    @{ Cell = 'D5'; Value = '0.07826' }
    @{ Cell = 'E5'; Value = '1.63%' }
    @{ Cell = 'B6'; Value = 'FTXToken' }
    @{ Cell = 'C6'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' }
    @{ Cell = 'D6'; Value = '1.904' }
    @{ Cell = 'E6'; Value = '0.66%' }
    @{ Cell = 'B7'; Value = 'KuCoinToken' }
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs' }
    @{ Cell = 'D7'; Value = '8.255' }
    @{ Cell = 'E7'; Value = '0.88%' }
    @{ Cell = 'D8'; Value = '2.952' }
    @{ Cell = 'E8'; Value = '-6.79%' }
    @{ Cell = 'D9'; Value = '0.9182' }
    @{ Cell = 'E9'; Value = '0.09%' }
    @{ Cell = 'D10'; Value = '0.1196' }
    @{ Cell = 'E10'; Value = '-0.17%' }
    @{ Cell = 'D11'; Value = '0.1911' }
    @{ Cell = 'E11'; Value = '3.23%' }
    @{ Cell = 'D12'; Value = '0.08959' }
    @{ Cell = 'E12'; Value = '3.91%' }
    @{ Cell = 'D13'; Value = '0.03352' }
    @{ Cell = 'E13'; Value = '-1.26%' }
    @{ Cell = 'D14'; Value = '0.09578' }
    @{ Cell = 'E14'; Value = '-1.30%' }
    @{ Cell = 'D15'; Value = '0.001380' }
    @{ Cell = 'E15'; Value = '0.27%' }
    @{ Cell = 'D16'; Value = '0.005704' }
    @{ Cell = 'E16'; Value = '-6.57%' }
    @{ Cell = 'D17'; Value = '3.534' }
    @{ Cell = 'E17'; Value = '-1.97%' }
    @{ Cell = 'D18'; Value = '4.422' }
    @{ Cell = 'E18'; Value = '1.62%' }
    @{ Cell = 'D19'; Value = '0.3440' }
    @{ Cell = 'E19'; Value = '0.92%' }
    @{ Cell = 'D20'; Value = '5.241' }
    @{ Cell = 'E20'; Value = '4.43%' }
    @{ Cell = 'D21'; Value = '0.1284' }
    @{ Cell = 'E21'; Value = '0.51%' }
    @{ Cell = 'D22'; Value = '0.2591' }
    @{ Cell = 'E22'; Value = '-0.35%' }
    @{ Cell = 'D23'; Value = '0.04361' }
    @{ Cell = 'E23'; Value = '0.63%' }
    @{ Cell = 'D24'; Value = '0.001250' }
    @{ Cell = 'E24'; Value = '2.99%' }
    @{ Cell = 'D25'; Value = '0.004659' }
    @{ Cell = 'E25'; Value = '10.89%' }
    @{ Cell = 'D26'; Value = '0.0001362' }
    @{ Cell = 'E26'; Value = '0.69%' }
    @{ Cell = 'D27'; Value = '0.0003991' }
    @{ Cell = 'E27'; Value = '-98.11%' }
    @{ Cell = 'D39'; Value = '0.02254' }
    @{ Cell = 'E39'; Value = '4.12%' }
    @{ Cell = 'D40'; Value = '0.05029' }
    @{ Cell = 'E40'; Value = '2.65%' }
    @{ Cell = 'D41'; Value = '0.007458' }
    @{ Cell = 'E41'; Value = '-1.08%' }
    @{ Cell = 'D42'; Value = '0.009045' }
    @{ Cell = 'E42'; Value = '-9.70%' }
    @{ Cell = 'E43'; Value = '1.19%' }
    @{ Cell = 'D44'; Value = '0.001953' }
    @{ Cell = 'E44'; Value = '-2.21%' }
    @{ Cell = 'D45'; Value = '0.009301' }
    @{ Cell = 'E45'; Value = '9.03%' }
    @{ Cell = 'D46'; Value = '0.00006578' }
    @{ Cell = 'E46'; Value = '0.21%' }
    @{ Cell = 'D47'; Value = '0.00000000750' }
    @{ Cell = 'E47'; Value = '-0.28%' }
    @{ Cell = 'D48'; Value = '0.003367' }
    @{ Cell = 'E48'; Value = '12.01%' }
    @{ Cell = 'D49'; Value = '0.001000' }
    @{ Cell = 'E49'; Value = '-23.30%' }
    @{ Cell = 'D50'; Value = '0.00002101' }
    @{ Cell = 'E50'; Value = '-0.28%' }
    @{ Cell = 'D51'; Value = '0.0002001' }
    @{ Cell = 'E51'; Value = '-0.28%' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = "'" + $u.Value
    $cell.Style = "Normal"
}
